# Actualizacion automatica 2025-11-26 13:30:09
# Insert a new client "CONSTRUCCION, INGENIERIA Y TECNOLOGIA CONSTRUINTEC SAS" into the
# "CASTRO ALCIVAR EDA MARIA" advisor block (alphabetically, right after
# "CONSTRUCCION MATUTE JIMENEZ CONSTRUMAJI S.A." at row 18) on both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, shifting all following client rows
# down by one. Also record a -1206.74 "NO RESURTIBLES" / "noviembre" adjustment for
# RUIZ PINEDA LUIS ALFREDO (now at row 49 after the shift), and refresh the
# "CUMPLIMIENTO MENSUAL" summary sheet accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row at position 18; existing rows 18..60 shift down to 19..61,
# and the footer row (old 61) shifts down to 62. Formatting of row 17 below is
# inherited automatically by the insert operation.
$ws1.Rows.Item(18).Insert()

$ws1.Range("A18").Value = "CASTRO ALCIVAR EDA MARIA"
$ws1.Range("B18").Value = "CONSTRUCCION, INGENIERIA Y TECNOLOGIA CONSTRUINTEC SAS"
$newRow1Cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $newRow1Cols) {
    $ws1.Range("${col}18").Value = 0
}

# Extra sales adjustment discovered for RUIZ PINEDA LUIS ALFREDO (now row 49,
# column P = "NO RESURTIBLES")
$ws1.Range("P49").Value = -1206.74

# Refresh the footer counts row (old row 61, now row 62): "X de 59" -> "X de 60"
$footerCols1 = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $footerCols1) {
    $cell = $ws1.Range("${col}62")
    $text = $cell.Value()
    $newText = $text -replace " de 59$", " de 60"
    $cell.Value = $newText
}

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(18).Insert()

$ws2.Range("A18").Value = "CASTRO ALCIVAR EDA MARIA"
$ws2.Range("B18").Value = "CONSTRUCCION, INGENIERIA Y TECNOLOGIA CONSTRUINTEC SAS"
$newRow2Cols = @("C","D","E","F","G")
foreach ($col in $newRow2Cols) {
    $ws2.Range("${col}18").Value = 0
}

# Same extra adjustment for RUIZ PINEDA LUIS ALFREDO (now row 49,
# column F = "noviembre")
$ws2.Range("F49").Value = -1206.74

# Refresh the TOTAL row (old row 61, now row 62)
$ws2.Range("C62").Value = 81440.64
$ws2.Range("D62").Value = 91039.07000000001
$ws2.Range("E62").Value = 91874.38
$ws2.Range("F62").Value = 56172.62
$ws2.Range("G62").Value = 85274.87999999999

# ---------------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 8 -> "NO RESURTIBLES" group totals
$ws3.Range("D8").Value = -1220.51
$ws3.Range("E8").Value = 1570.51
$ws3.Range("F8").Value = -3.487171428571429

# Row 14 -> TOTAL
$ws3.Range("D14").Value = 59138.56
$ws3.Range("E14").Value = 38723.32766749098
$ws3.Range("F14").Value = 0.6043063485647988

# Column F width adjusts slightly (autofit after the VENTA/CUMPLIMIENTO values changed)
$ws3.Columns.Item(6).ColumnWidth = 24
